$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for changed rows.
# A leading apostrophe forces Excel to keep values that look like plain
# numbers (e.g. "227.87") stored as text, matching the source data which
# uses inline/shared strings for all Price cells.
$ws.Range("D2").Value = "39.661.68"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").Value = "2.162.96"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'227.87"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'0.632"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("D7").Value = "'63.51"
$ws.Range("E7").Value = "  +1.61%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.394"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").Value = "'0.0850"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "'16.01"
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("D13").Value = "2.484.60"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").Value = "'22.00"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "'0.809"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "'5.49"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "2.168.91"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "39.627.32"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").Value = "'71.94"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "'6.13"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "0.0₃0847"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  +0.04%  "

# Rows 24 and 25 swap: PancakeSwap moves to rank 24, Toncoin moves to rank 25
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.40"
$ws.Range("E24").Value = "  +3.47%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("E25").Value = "  +0.91%  "

$ws.Range("D26").Value = "'9.65"
$ws.Range("D27").Value = "'172.23"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D29").Value = "'19.80"
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "'2.68"
$ws.Range("E31").Value = "  +3.39%  "
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").Value = "'4.70"
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").Value = "'0.0619"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("D38").Value = "'3.63"
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'4.76"
$ws.Range("E40").Value = "  +13.27%  "
$ws.Range("D41").Value = "'102.18"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "'17.66"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("D44").Value = "1.511.14"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("E45").Value = "  +1.30%  "
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("E48").Value = "  +1.94%  "
$ws.Range("D49").Value = "'7.75"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("E50").Value = "  +1.19%  "

# Row 51: TerraClassic replaced by RocketPoolETH
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.369.18"
$ws.Range("E51").Value = "  +2.86%  "
